{"js": "// Apply two paragraph-text edits described by the commit:\n//  1. \"Reservations can be cancelled only two hours before the ride.\"\n//     -> \"Requests and reservations can be cancelled if and only if no\n//         taxi have been assigned to the customer.\"\n//  2. \"...for customers only. Taxi drivers' account are created when they\n//      are hired by the taxi company. Taxi drivers will receive...\"\n//     -> \"...for customers only. Taxi drivers' account are created by an\n//         administrator when they are hired by the taxi company. Taxi\n//         drivers will then receive...\"\n//\n// Because the target shape of each paragraph needs specific run breaks and\n// w:proofErr grammar-check markers (which are not modeled by the Office.js\n// object model), we rebuild each paragraph's contents with insertOoxml so\n// the resulting run/proofErr structure matches exactly, while leaving the\n// paragraph's own identity (and every other paragraph) untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst RFONTS = '<w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr>';\n\nfunction wrapPkg(innerParagraphXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" + innerParagraphXml + \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// --- Paragraph 1: cancellation-rule bullet ---------------------------------\nconst cancelPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Reservations can be cancelled only two hours before the ride.\") !== -1\n);\nif (!cancelPara) {\n  throw new Error(\"Could not find the 'Reservations can be cancelled...' paragraph\");\n}\n\nconst cancelOoxml = wrapPkg(\n  '<w:p w:rsidR=\"004934BF\" w:rsidRDefault=\"009D6D77\" w:rsidP=\"009D6D77\">' +\n    '<w:pPr><w:pStyle w:val=\"Paragrafoelenco\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    RFONTS +\n    \"</w:pPr>\" +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\">Requests and reservations </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r>\" + RFONTS + \"<w:t>can be cancelled</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\"> if and only if no taxi have been assigned to the customer.</w:t></w:r>' +\n    \"</w:p>\"\n);\n\ncancelPara.insertOoxml(cancelOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Paragraph 2: taxi-driver account creation bullet -----------------------\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst accountsPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Taxi drivers\\u2019 account are created\") !== -1\n);\nif (!accountsPara) {\n  throw new Error(\"Could not find the 'Taxi drivers' account are created...' paragraph\");\n}\n\nconst accountsOoxml = wrapPkg(\n  '<w:p w:rsidR=\"00B75ECE\" w:rsidRDefault=\"00B018D6\" w:rsidP=\"00B75ECE\">' +\n    '<w:pPr><w:pStyle w:val=\"Paragrafoelenco\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr>' +\n    RFONTS +\n    \"</w:pPr>\" +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\">The web and mobile registration </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r>\" + RFONTS + \"<w:t>is intended</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\"> for customers only. </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    \"<w:r>\" + RFONTS + '<w:t>Taxi drivers\\u2019 account are created</w:t></w:r>' +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\"> by an administrator</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\"> when they are hired by the taxi company. Taxi drivers will</w:t></w:r>' +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\"> then</w:t></w:r>' +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\"> receive </w:t></w:r>' +\n    '<w:r w:rsidR=\"00E55C4B\">' + RFONTS + \"<w:t>their</w:t></w:r>\" +\n    \"<w:r>\" + RFONTS + '<w:t xml:space=\"preserve\"> username and temporary password, which </w:t></w:r>' +\n    '<w:r w:rsidR=\"00E55C4B\">' + RFONTS + '<w:t xml:space=\"preserve\">they will be able to change </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r w:rsidR=\"00E55C4B\">' + RFONTS + \"<w:t>once</w:t></w:r>\" +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r w:rsidR=\"00E55C4B\">' + RFONTS + '<w:t xml:space=\"preserve\"> logged in the application.</w:t></w:r>' +\n    \"</w:p>\"\n);\n\naccountsPara.insertOoxml(accountsOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Apply two paragraph-text edits described by the commit:\n#  1. \"Reservations can be cancelled only two hours before the ride.\"\n#     -> \"Requests and reservations can be cancelled if and only if no\n#         taxi have been assigned to the customer.\"\n#  2. \"...for customers only. Taxi drivers' account are created when they\n#      are hired by the taxi company. Taxi drivers will receive...\"\n#     -> \"...for customers only. Taxi drivers' account are created by an\n#         administrator when they are hired by the taxi company. Taxi\n#         drivers will then receive...\"\n#\n# Each paragraph's body is rebuilt with Range.InsertXML so the resulting\n# run / w:proofErr (grammar-check marker) structure matches the target\n# OOXML exactly; InsertXML replaces only the addressed range's contents,\n# leaving every other paragraph untouched.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphContaining([string]$needle) {\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.Contains($needle)) {\n            return $p\n        }\n    }\n    throw \"Paragraph containing '$needle' not found\"\n}\n\n# --- Paragraph: cancellation-rule bullet -----------------------------------\n$cancelPara = Find-ParagraphContaining \"Reservations can be cancelled only two hours before the ride.\"\n\n$cancelXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n'<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n'<w:p w:rsidR=\"004934BF\" w:rsidRDefault=\"009D6D77\" w:rsidP=\"009D6D77\">' +\n'<w:pPr><w:pStyle w:val=\"Paragrafoelenco\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr></w:pPr>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\">Requests and reservations </w:t></w:r>' +\n'<w:proofErr w:type=\"gramStart\"/>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t>can be cancelled</w:t></w:r>' +\n'<w:proofErr w:type=\"gramEnd\"/>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\"> if and only if no taxi have been assigned to the customer.</w:t></w:r>' +\n'</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$cancelPara.Range.InsertXML($cancelXml)\n\n# --- Paragraph: taxi-driver account creation bullet -------------------------\n$accountsPara = Find-ParagraphContaining \"Taxi drivers\u2019 account are created\"\n\n$accountsXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n'<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n'<w:p w:rsidR=\"00B75ECE\" w:rsidRDefault=\"00B018D6\" w:rsidP=\"00B75ECE\">' +\n'<w:pPr><w:pStyle w:val=\"Paragrafoelenco\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr></w:pPr>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\">The web and mobile registration </w:t></w:r>' +\n'<w:proofErr w:type=\"gramStart\"/>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t>is intended</w:t></w:r>' +\n'<w:proofErr w:type=\"gramEnd\"/>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\"> for customers only. </w:t></w:r>' +\n'<w:proofErr w:type=\"gramStart\"/>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t>Taxi drivers\u2019 account are created</w:t></w:r>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\"> by an administrator</w:t></w:r>' +\n'<w:proofErr w:type=\"gramEnd\"/>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\"> when they are hired by the taxi company. Taxi drivers will</w:t></w:r>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\"> then</w:t></w:r>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\"> receive </w:t></w:r>' +\n'<w:r w:rsidR=\"00E55C4B\"><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t>their</w:t></w:r>' +\n'<w:r><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\"> username and temporary password, which </w:t></w:r>' +\n'<w:r w:rsidR=\"00E55C4B\"><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\">they will be able to change </w:t></w:r>' +\n'<w:proofErr w:type=\"gramStart\"/>' +\n'<w:r w:rsidR=\"00E55C4B\"><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t>once</w:t></w:r>' +\n'<w:proofErr w:type=\"gramEnd\"/>' +\n'<w:r w:rsidR=\"00E55C4B\"><w:rPr><w:rFonts w:ascii=\"Georgia\" w:hAnsi=\"Georgia\"/></w:rPr><w:t xml:space=\"preserve\"> logged in the application.</w:t></w:r>' +\n'</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$accountsPara.Range.InsertXML($accountsXml)\n"}
